$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sector>Financiamiento e información financiera de intermediarios financieros>Cuadros estadísticos descontinuados>Balances financieros sectoriales>Banca comercial (metodología 2000)>Recursos y obligaciones>Banca Interna>9.2.1 De Intermediarios Financieros>Recursos de la banca interna Saldos nominales en miles de pesos Recursos totales 9.2 Activos de larga duración disponibles para la venta 9.2.1 De intermediarios financieros"

$ws.Range("A3").Value = "Sector>Financiamiento e información financiera de intermediarios financieros>Financiamiento>Financiamiento otorgado por la banca comercial>Crédito por la principal actividad del acreditado (total)>Moneda nacional>Comercio al por mayor de productos farmacéuticos, de perfumería, artículos para el esparcimiento, electrodomésticos meno"

$ws.Range("A4").Value = "Sector>Financiamiento e información financiera de intermediarios financieros>Captación>Banca comercial>por entidad federativa y municipio>saldos de los principales instrumentos de captación>Cuautitlán Izcalli>Saldos de los instrumentos de captación de la Banca Comercial por entidad federativa y municipio 1/ *Instrumentos de captación *Saldos de captación 2/ Memo: Cuentas con saldos a favor de Tarjeta de Crédito 5/ México Cuautitlán Izcalli"

$ws.Range("A5").Value = "Sector>Financiamiento e información financiera de intermediarios financieros>Cuadros estadísticos descontinuados>Captación>Banca comercial>Regional>Agencias en el exterior (metodología 2008)>Número de Oficinas"

$ws.Range("A6").Value = "Sector>Financiamiento e información financiera de intermediarios financieros>Cuadros estadísticos descontinuados>Captación>Banca de desarrollo>Regional>Total, republica y agencias (metodología 2008)>Guanajuato>Pagarés con rendimiento liquidable al vencimiento>Silao *Captación tradicional *Saldos de captación en moneda extranjera 2/ Pagarés con rendimiento liquidable al vencimiento"
